# Add season record columns (Wins, Losses, Ties) to the team statistics sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, centered, bordered) from AC1
# onto the new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-48): same Wins/Losses/Ties record for every player row ---
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 95   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 67   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
